$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 241.65384
$ws.Range("I28").Value = 205.5
$ws.Range("J28").Value = 323
$ws.Range("K28").Value = 205.5
$ws.Range("L28").Value = 323
$ws.Range("M28").Value = 279.5
$ws.Range("N28").Value = -1293

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3143.6
$ws.Range("I64").Value = 3185.8572
$ws.Range("J64").Value = 3106.625
$ws.Range("K64").Value = 3185.8572
$ws.Range("L64").Value = 3106.625
$ws.Range("M64").Value = -2937.8572
$ws.Range("N64").Value = -3602.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3143.6
$ws.Range("I67").Value = 3185.8572
$ws.Range("J67").Value = 3106.625
$ws.Range("K67").Value = 3185.8572
$ws.Range("L67").Value = 3106.625
$ws.Range("M67").Value = -2327.8572
$ws.Range("N67").Value = -4822.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3692.158
$ws.Range("I69").Value = 5200.3335
$ws.Range("J69").Value = 3409.375
$ws.Range("K69").Value = 15601.0005
$ws.Range("L69").Value = 10228.125
$ws.Range("M69").Value = -14727.0005
$ws.Range("N69").Value = -11976.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 3692.158
$ws.Range("I72").Value = 5200.3335
$ws.Range("J72").Value = 3409.375
$ws.Range("K72").Value = 46803.0015
$ws.Range("L72").Value = 30684.375
$ws.Range("M72").Value = -42435.0015
$ws.Range("N72").Value = -39420.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4288.8887
$ws.Range("I76").Value = 3500
$ws.Range("J76").Value = 4514.2856
$ws.Range("K76").Value = 3500
$ws.Range("L76").Value = 4514.2856
$ws.Range("M76").Value = -3185
$ws.Range("N76").Value = -5144.2856

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4288.8887
$ws.Range("I79").Value = 3500
$ws.Range("J79").Value = 4514.2856
$ws.Range("K79").Value = 3500
$ws.Range("L79").Value = 4514.2856
$ws.Range("M79").Value = -2408
$ws.Range("N79").Value = -6698.2856

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 9259946
$ws.Range("I92").Value = 20833640
$ws.Range("J92").Value = 990.5
$ws.Range("K92").Value = 20833640
$ws.Range("L92").Value = 990.5
$ws.Range("M92").Value = -20832392

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 30601
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 30601
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 30601
$ws.Range("N93").Value = -35593

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 5000
$ws.Range("I98").Value = 5000
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 5000
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -3502
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 33333580
$ws.Range("I107").Value = 47619228
$ws.Range("J107").Value = 404.44446
$ws.Range("K107").Value = 47619228
$ws.Range("L107").Value = 404.44446
$ws.Range("M107").Value = -47617308
$ws.Range("N107").Value = -4244.44446

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H114").Value = 98000
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 98000
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 98000
$ws.Range("N114").Value = -106678

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H117").Value = 58000
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 58000
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 58000
$ws.Range("N117").Value = -67178

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -12550
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1117.317
$ws.Range("I137").Value = 850.5294
$ws.Range("J137").Value = 1556.2258
$ws.Range("K137").Value = 2551.5882
$ws.Range("L137").Value = 4668.6774
$ws.Range("M137").Value = -1.588200000000143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1051.5
$ws.Range("I2").Value = 1051.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1051.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -938.5
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 507.5
$ws.Range("I14").Value = 507.5
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 507.5
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -332.5
$ws.Range("N14").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 27500
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 27500
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 27500
$ws.Range("N24").Value = -28248

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1754.4445
$ws.Range("I45").Value = 1700.6
$ws.Range("J45").Value = 1821.75
$ws.Range("K45").Value = 1700.6
$ws.Range("L45").Value = 1821.75
$ws.Range("M45").Value = -1323.6
$ws.Range("N45").Value = -2575.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2288
$ws.Range("I61").Value = 1419.0588
$ws.Range("J61").Value = 4750
$ws.Range("K61").Value = 1419.0588
$ws.Range("L61").Value = 4750
$ws.Range("M61").Value = -1207.0588
$ws.Range("N61").Value = -5174

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 23720840
$ws.Range("I74").Value = 25002964
$ws.Range("J74").Value = 20836060
$ws.Range("K74").Value = 25002964
$ws.Range("L74").Value = 20836060
$ws.Range("M74").Value = -25002090
$ws.Range("N74").Value = -20837808

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 23720840
$ws.Range("I77").Value = 25002964
$ws.Range("J77").Value = 20836060
$ws.Range("K77").Value = 125014820
$ws.Range("L77").Value = 104180300
$ws.Range("M77").Value = -125010452
$ws.Range("N77").Value = -104189036

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 26736.6
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 26736.6
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 26736.6
$ws.Range("N96").Value = -32228.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H100").Value = 27500
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 27500
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 27500
$ws.Range("N100").Value = -29664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 5699.778
$ws.Range("I110").Value = 5953.9165
$ws.Range("J110").Value = 3666.6667
$ws.Range("K110").Value = 5953.9165
$ws.Range("L110").Value = 3666.6667
$ws.Range("M110").Value = -3908.9165
$ws.Range("N110").Value = -7756.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1051.5
$ws.Range("I116").Value = 1051.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1051.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1242.5
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2288
$ws.Range("I136").Value = 1419.0588
$ws.Range("J136").Value = 4750
$ws.Range("K136").Value = 4257.1764
$ws.Range("L136").Value = 14250
$ws.Range("M136").Value = -1707.1764
$ws.Range("N136").Value = -19350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1051.5
$ws.Range("I3").Value = 1051.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1051.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -937.5
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 122.57143
$ws.Range("I80").Value = 132
$ws.Range("J80").Value = 113.14286
$ws.Range("K80").Value = 132
$ws.Range("L80").Value = 113.14286
$ws.Range("M80").Value = 866
$ws.Range("N80").Value = -2109.14286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 122.57143
$ws.Range("I83").Value = 132
$ws.Range("J83").Value = 113.14286
$ws.Range("K83").Value = 660
$ws.Range("L83").Value = 565.7143
$ws.Range("M83").Value = 4332
$ws.Range("N83").Value = -10549.7143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3477.375
$ws.Range("I107").Value = 2990
$ws.Range("J107").Value = 3639.8333
$ws.Range("K107").Value = 2990
$ws.Range("L107").Value = 3639.8333
$ws.Range("M107").Value = -1070
$ws.Range("N107").Value = -7479.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7845075.5
$ws.Range("I31").Value = 9092156
$ws.Range("J31").Value = 6899014
$ws.Range("K31").Value = 9092156
$ws.Range("L31").Value = 6899014
$ws.Range("M31").Value = -9091861
$ws.Range("N31").Value = -6899604

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7845075.5
$ws.Range("I34").Value = 9092156
$ws.Range("J34").Value = 6899014
$ws.Range("K34").Value = 9092156
$ws.Range("L34").Value = 6899014
$ws.Range("M34").Value = -9091954
$ws.Range("N34").Value = -6899418

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1193923.1
$ws.Range("I99").Value = 1556010.6
$ws.Range("J99").Value = 4207.143
$ws.Range("K99").Value = 1556010.6
$ws.Range("L99").Value = 4207.143
$ws.Range("M99").Value = -1554512.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1193923.1
$ws.Range("I126").Value = 1556010.6
$ws.Range("J126").Value = 4207.143
$ws.Range("K126").Value = 4668031.800000001
$ws.Range("L126").Value = 12621.429
$ws.Range("M126").Value = -4665561.800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2511.375
$ws.Range("I132").Value = 1128.4445
$ws.Range("J132").Value = 4289.4287
$ws.Range("K132").Value = 3385.3335
$ws.Range("L132").Value = 12868.2861
$ws.Range("M132").Value = -855.3335000000002
$ws.Range("N132").Value = -17928.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1191359.5
$ws.Range("I131").Value = 3509452.8
$ws.Range("J131").Value = 987.1622
$ws.Range("K131").Value = 10528358.4
$ws.Range("L131").Value = 2961.4866
$ws.Range("M131").Value = -10523318.4
$ws.Range("N131").Value = -13041.4866

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1652.4
$ws.Range("I113").Value = 1168.4117
$ws.Range("J113").Value = 2010.1305
$ws.Range("K113").Value = 1168.4117
$ws.Range("L113").Value = 2010.1305
$ws.Range("M113").Value = 1001.5883
$ws.Range("N113").Value = -6350.1305

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 994.5
$ws.Range("I122").Value = 994.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2983.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -533.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4870.9
$ws.Range("I40").Value = 4838.625
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 4838.625
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -4702.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 20349
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 20349
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 20349
$ws.Range("N119").Value = -30025
